# Scheduled-runner update: refresh cached Universalis market-price snapshots
# and recomputed leve-profit figures across the per-class profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 404.33334
$ws.Range("I33").Value = 404.33334
$ws.Range("K33").Value = 404.33334
$ws.Range("M33").Value = -175.33334
$ws.Range("H137").Value = 3217.3
$ws.Range("I137").Value = 1366
$ws.Range("K137").Value = 4098
$ws.Range("M137").Value = -1548

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2709.3333
$ws.Range("I32").Value = 2709.3333
$ws.Range("K32").Value = 2709.3333
$ws.Range("M32").Value = -2422.3333
$ws.Range("H74").Value = 2754.2222
$ws.Range("I74").Value = 688.7273
$ws.Range("K74").Value = 688.7273
$ws.Range("M74").Value = 185.2727
$ws.Range("H77").Value = 2754.2222
$ws.Range("I77").Value = 688.7273
$ws.Range("K77").Value = 3443.6365
$ws.Range("M77").Value = 924.3634999999999
$ws.Range("H97").Value = 862.8889
$ws.Range("I97").Value = 696.8570999999999
$ws.Range("K97").Value = 696.8570999999999
$ws.Range("M97").Value = -200.8570999999999
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("H122").Value = 3355.2942
$ws.Range("I122").Value = 3138.8
$ws.Range("K122").Value = 9416.400000000001
$ws.Range("M122").Value = -6966.400000000001
$ws.Range("H132").Value = 3524.375
$ws.Range("I132").Value = 2609.2222
$ws.Range("K132").Value = 7827.6666
$ws.Range("M132").Value = -5297.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3525.5
$ws.Range("I20").Value = 3936.125
$ws.Range("K20").Value = 3936.125
$ws.Range("M20").Value = -3689.125
$ws.Range("H86").Value = 2651.1667
$ws.Range("I86").Value = 2481.4
$ws.Range("J86").Value = 3500
$ws.Range("K86").Value = 2481.4
$ws.Range("L86").Value = 3500
$ws.Range("M86").Value = -1358.4
$ws.Range("N86").Value = -5746
$ws.Range("H89").Value = 2651.1667
$ws.Range("I89").Value = 2481.4
$ws.Range("J89").Value = 3500
$ws.Range("K89").Value = 12407
$ws.Range("L89").Value = 17500
$ws.Range("M89").Value = -6791
$ws.Range("N89").Value = -28732
$ws.Range("H94").Value = 1439.4286
$ws.Range("I94").Value = 1519
$ws.Range("J94").Value = 1333.3334
$ws.Range("K94").Value = 1519
$ws.Range("L94").Value = 1333.3334
$ws.Range("M94").Value = -1068
$ws.Range("N94").Value = -2235.3334
$ws.Range("H106").Value = 5494
$ws.Range("J106").Value = 5494
$ws.Range("L106").Value = 5494
$ws.Range("N106").Value = -8018
$ws.Range("H134").Value = 2567.2727
$ws.Range("I134").Value = 2437.8333
$ws.Range("J134").Value = 3149.75
$ws.Range("K134").Value = 7313.499899999999
$ws.Range("L134").Value = 9449.25
$ws.Range("M134").Value = -4778.499899999999
$ws.Range("N134").Value = -14519.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5658.3335
$ws.Range("I16").Value = 6987.5
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 6987.5
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -6700.5
$ws.Range("N16").Value = -3574
$ws.Range("H31").Value = 2439.6316
$ws.Range("I31").Value = 1024.125
$ws.Range("J31").Value = 9989
$ws.Range("K31").Value = 1024.125
$ws.Range("L31").Value = 9989
$ws.Range("M31").Value = -729.125
$ws.Range("N31").Value = -10579
$ws.Range("H34").Value = 2439.6316
$ws.Range("I34").Value = 1024.125
$ws.Range("J34").Value = 9989
$ws.Range("K34").Value = 1024.125
$ws.Range("L34").Value = 9989
$ws.Range("M34").Value = -822.125
$ws.Range("N34").Value = -10393
$ws.Range("H58").Value = 1309
$ws.Range("I58").Value = 1368.7778
$ws.Range("K58").Value = 1368.7778
$ws.Range("M58").Value = -1165.7778
$ws.Range("H113").Value = 5658.3335
$ws.Range("I113").Value = 6987.5
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 6987.5
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -4817.5
$ws.Range("N113").Value = -7340
$ws.Range("H134").Value = 3984.0833
$ws.Range("I134").Value = 3996.6
$ws.Range("J134").Value = 3921.5
$ws.Range("K134").Value = 11989.8
$ws.Range("L134").Value = 11764.5
$ws.Range("M134").Value = -9454.799999999999
$ws.Range("N134").Value = -16834.5
$ws.Range("H136").Value = 1309
$ws.Range("I136").Value = 1368.7778
$ws.Range("K136").Value = 4106.3334
$ws.Range("M136").Value = -1556.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 200000
$ws.Range("J37").Value = 200000
$ws.Range("L37").Value = 600000
$ws.Range("N37").Value = -600224
$ws.Range("H109").Value = 675.3333
$ws.Range("I109").Value = 675.3333
$ws.Range("K109").Value = 2025.9999
$ws.Range("M109").Value = -985.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 66.25
$ws.Range("I2").Value = 71.666664
$ws.Range("K2").Value = 71.666664
$ws.Range("M2").Value = 41.333336
$ws.Range("H25").Value = 2000
$ws.Range("I25").Value = 2000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 2000
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -1471
$ws.Range("N25").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -705
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 1000
$ws.Range("M27").Value = -893
$ws.Range("H68").Value = 2995.1
$ws.Range("I68").Value = 2995.1
$ws.Range("K68").Value = 2995.1
$ws.Range("M68").Value = -2246.1
$ws.Range("H71").Value = 2995.1
$ws.Range("I71").Value = 2995.1
$ws.Range("K71").Value = 14975.5
$ws.Range("M71").Value = -11231.5
$ws.Range("H93").Value = 1186.2
$ws.Range("I93").Value = 978.6667
$ws.Range("J93").Value = 1497.5
$ws.Range("K93").Value = 978.6667
$ws.Range("L93").Value = 1497.5
$ws.Range("M93").Value = 269.3333
$ws.Range("N93").Value = -3993.5
$ws.Range("H122").Value = 3001.3333
$ws.Range("I122").Value = 3001.3333
$ws.Range("K122").Value = 9003.999899999999
$ws.Range("M122").Value = -6553.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 372
$ws.Range("I107").Value = 196.33333
$ws.Range("J107").Value = 899
$ws.Range("K107").Value = 588.99999
$ws.Range("L107").Value = 2697
$ws.Range("M107").Value = 1331.00001
$ws.Range("N107").Value = -6537
